# Update "want to go" counts / price values across the sheets to reflect
# the refreshed scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet - column F updates
$wsExhibit.Range("F5").Value  = 1047
$wsExhibit.Range("F13").Value = 526
$wsExhibit.Range("F14").Value = 1693
$wsExhibit.Range("F15").Value = 1401
$wsExhibit.Range("F16").Value = 822
$wsExhibit.Range("F21").Value = 1161
$wsExhibit.Range("F22").Value = 144
$wsExhibit.Range("F24").Value = 40
$wsExhibit.Range("F25").Value = 3603
$wsExhibit.Range("F29").Value = 53

# 演出 (Show) sheet - column F and G updates
$wsShow.Range("F8").Value  = 36
$wsShow.Range("F9").Value  = 32
$wsShow.Range("G13").Value = 180

# 全部类型 (All types) sheet - column F and G updates
$wsAll.Range("F13").Value = 36
$wsAll.Range("F14").Value = 32
$wsAll.Range("F16").Value = 1047
$wsAll.Range("F24").Value = 526
$wsAll.Range("F25").Value = 1693
$wsAll.Range("F26").Value = 1401
$wsAll.Range("F27").Value = 822
$wsAll.Range("F34").Value = 1161
$wsAll.Range("F35").Value = 144
$wsAll.Range("F37").Value = 40
$wsAll.Range("F38").Value = 3603
$wsAll.Range("G42").Value = 180
$wsAll.Range("F44").Value = 53
